# Apply the commit's data changes to Sheet1:
#  - Row 2 (A2:B2) changes from mngr195759/AhuzYte -> mngr201383/jAzaryp
#  - A new Row 5 (A5:B5) is inserted with mngr201383/jAzaryp
#  - The old Row 5 becomes Row 6, keeping A6 = mngr1957 but B6 changes
#    from AhuzYtek -> AhuzYte
#  - Rows 1, 3, 4 are unchanged
#  - Selection ends on C6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 5 so the old row5 data shifts down to row6
$ws.Rows.Item(5).Insert()

# Row 2: update credentials
$ws.Range("A2").Value = "mngr201383"
$ws.Range("B2").Value = "jAzaryp"

# New Row 5: same credentials as the updated row 2
$ws.Range("A5").Value = "mngr201383"
$ws.Range("B5").Value = "jAzaryp"

# Row 6 (previously row 5, shifted down by the insert): fix up B6
$ws.Range("A6").Value = "mngr1957"
$ws.Range("B6").Value = "AhuzYte"

# Match the final selection recorded in the workbook
$ws.Range("C6").Select()
